# Updated cryptos list on Wed Nov 27 18:00:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces Excel to store these as text (matching the original
# inlineStr cells) instead of auto-coercing numeric-looking strings (e.g.
# "240.08", "1.49") into actual numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'96.257.27"
$ws.Range("E2").Value = "'  +2.84%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.562.18"
$ws.Range("E3").Value = "'  +7.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "'  -0.03%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'240.08"
$ws.Range("E5").Value = "'  +4.18%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'637.53"
$ws.Range("E6").Value = "'  +3.22%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'1.49"
$ws.Range("E7").Value = "'  +7.68%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.402"
$ws.Range("E8").Value = "'  +3.80%  "

# Row 9 - USDC
$ws.Range("E9").Value = "'  -0.05%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "'  +10.41%  "

# Row 11 - LidoStakedEther
$ws.Range("D11").Value = "'3.560.83"
$ws.Range("E11").Value = "'  +6.93%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'43.33"
$ws.Range("E12").Value = "'  +2.44%  "

# Row 13 - TRON
$ws.Range("E13").Value = "'  +4.07%  "

# Row 14 - Toncoin
$ws.Range("E14").Value = "'  +8.38%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.227.41"
$ws.Range("E15").Value = "'  +7.04%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'96.152.04"
$ws.Range("E16").Value = "'  +2.82%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "'  +4.41%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.564.22"
$ws.Range("E18").Value = "'  +7.02%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "'  +20.57%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'8.01"
$ws.Range("E20").Value = "'  -0.56%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'18.11"
$ws.Range("E21").Value = "'  +5.39%  "

# Row 22 - Stellar
$ws.Range("D22").Value = "'0.507"
$ws.Range("E22").Value = "'  +13.81%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'516.24"
$ws.Range("E23").Value = "'  +4.51%  "

# Row 24 - SuiNetwork
$ws.Range("D24").Value = "'3.44"
$ws.Range("E24").Value = "'  -0.12%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "'  +7.59%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "'  +8.88%  "

# Row 27 - Litecoin
$ws.Range("D27").Value = "'96.94"
$ws.Range("E27").Value = "'  +6.08%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "'12.36"
$ws.Range("E28").Value = "'  +5.59%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "'3.11"
$ws.Range("E29").Value = "'  +19.05%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "'  +5.22%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'11.56"
$ws.Range("E31").Value = "'  +4.70%  "

# Row 32 - Dai
$ws.Range("E32").Value = "'  -0.02%  "

# Row 33 - Cronos
$ws.Range("E33").Value = "'  +5.35%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "'  -0.43%  "

# Row 35 - EthereumClassic
$ws.Range("D35").Value = "'30.19"
$ws.Range("E35").Value = "'  +6.76%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "'  +6.90%  "

# Row 37 - Bittensor
$ws.Range("D37").Value = "'580.42"
$ws.Range("E37").Value = "'  +8.71%  "

# Row 38 & 39 - swap RenderToken / Fetch.AI order plus new values
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.51"
$ws.Range("E38").Value = "'  +10.78%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.91"
$ws.Range("E39").Value = "'  +7.32%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "'  +2.83%  "

# Row 41 - USDe
$ws.Range("E41").Value = "'  -0.02%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "'  +7.07%  "

# Row 43 & 44 - swap VeChain / ImmutableX order plus new values
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "'1.76"
$ws.Range("E43").Value = "'  +4.73%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0434"
$ws.Range("E44").Value = "'  +5.17%  "

# Row 45 - WhiteBITCoin
$ws.Range("E45").Value = "'  -0.79%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "'  +4.41%  "

# Row 47 - MantraDAO
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "'  -3.38%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "'  +4.03%  "

# Row 49 - OKB
$ws.Range("E49").Value = "'  +3.50%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "'8.20"
$ws.Range("E50").Value = "'  +3.20%  "

# Row 51 - dogwifhat
$ws.Range("E51").Value = "'  +2.96%  "
